# Append the new broker/quarter/ticker FVTPL-value rows (rows 44-67) that
# were added to Sheet1, and leave the selection on the last-edited cell
# (A56) to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (Broker, Quarter, Ticker, FVTPL value)
$data = @(
    ,("CTS", "1Q25", "VSC",    189.8)
    ,("CTS", "1Q25", "GEX",    106.4)
    ,("CTS", "1Q25", "EIB",    204.8)
    ,("CTS", "1Q25", "PET",    79.1)
    ,("CTS", "1Q25", "VPB",    70.6)
    ,("CTS", "1Q25", "Others", 0.1)
    ,("CTS", "4Q24", "VSC",    170.9)
    ,("CTS", "4Q24", "EIB",    238.4)
    ,("CTS", "4Q24", "GEX",    123)
    ,("CTS", "4Q24", "PET",    91.6)
    ,("CTS", "4Q24", "VPB",    71.3)
    ,("CTS", "4Q24", "FCN",    47.6)
    ,("CTS", "4Q24", "TTC",    178.9)
    ,("CTS", "4Q24", "PLC",    10.3)
    ,("VIX", "4Q24", "EIB",    1635)
    ,("VIX", "4Q24", "GEX",    787)
    ,("VIX", "4Q24", "VSC",    734.97)
    ,("VIX", "4Q24", "NVL",    0)
    ,("VIX", "4Q24", "HAH",    683.3)
    ,("VIX", "4Q24", "GEE",    428.9)
    ,("VIX", "4Q24", "Others", 4362.449)
    ,("VIX", "4Q24", "BSR",    684.672)
    ,("VIX", "4Q24", "SEA",    873)
    ,("VIX", "4Q24", "GEI",    298.214)
)

$startRow = 44
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Match the saved cursor/selection position from the edited workbook.
[void]$ws.Range("A56").Select()

Write-Output "Added $($data.Count) rows (44-$($startRow + $data.Count - 1)) to Sheet1"
